$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "MuSCs"-as-sender rows (rows 6-9) entirely; the new TPM
# re-run only produced 4 data rows (2-5), all sent from FAPs.
$ws.Range("A6:T9").EntireRow.Delete()

# Refresh the recomputed NATMI TPM metrics for the 4 remaining rows.
# Row 2: FAPs -> Inha -> Acvr2a -> ECs
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 5.566863615928334
$ws.Range("R2").Value = 50.101772543355
$ws.Range("S2").Value = 0.2087950866344732
$ws.Range("T2").Value = 0.2087950866344732

# Row 3: FAPs -> Inha -> Acvr2a -> FAPs
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("S3").Value = 0.4253229592313036
$ws.Range("T3").Value = 0.4253229592313036

# Row 4: FAPs -> Inha -> Acvr2a -> MuSCs
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 7.852380505081666
$ws.Range("R4").Value = 70.67142454573499
$ws.Range("S4").Value = 0.2945174484164121
$ws.Range("T4").Value = 0.2945174484164122

# Row 5: FAPs -> Inha -> Acvr2a -> Resolving-Mac
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 1.90270986139
$ws.Range("R5").Value = 17.12438875251
$ws.Range("S5").Value = 0.07136450571781097
$ws.Range("T5").Value = 0.07136450571781099
